# Refresh the crypto price/volume table with the latest scraped values.
# Cells D/B/C/E hold plain text (coin names, coinranking.com URLs, price
# strings such as '29.417.26' or '1.0000', and padded percentage strings).
# A leading apostrophe forces Excel to keep number-looking price strings as
# text instead of silently converting them to numeric values (which would
# drop meaningful trailing zeros, e.g. '4.130' -> 4.13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.417.26'
$ws.Range("E2").Value = '  -0.20%  '
# Row 3
$ws.Range("D3").Value = '1.849.50'
$ws.Range("E3").Value = '  -0.01%  '
# Row 4
$ws.Range("D4").Value = '''0.9989'
$ws.Range("E4").Value = '  -0.09%  '
# Row 5
$ws.Range("D5").Value = '''240.90'
$ws.Range("E5").Value = '  -0.93%  '
# Row 6
$ws.Range("D6").Value = '''0.6329'
$ws.Range("E6").Value = '  -1.61%  '
# Row 7
$ws.Range("D7").Value = '''1.0000'
$ws.Range("E7").Value = '  -0.02%  '
# Row 8
$ws.Range("D8").Value = '3.484.90'
$ws.Range("E8").Value = '  +88.16%  '
# Row 9
$ws.Range("D9").Value = '''0.07555'
$ws.Range("E9").Value = '  +1.34%  '
# Row 10
$ws.Range("D10").Value = '''0.2971'
$ws.Range("E10").Value = '  -0.58%  '
# Row 11
$ws.Range("D11").Value = '''24.66'
$ws.Range("E11").Value = '  +1.75%  '
# Row 12
$ws.Range("D12").Value = '3.670.38'
$ws.Range("E12").Value = '  +77.41%  '
# Row 13
$ws.Range("D13").Value = '''0.07711'
$ws.Range("E13").Value = '  +1.08%  '
# Row 14
$ws.Range("D14").Value = '''4.987'
$ws.Range("E14").Value = '  -0.57%  '
# Row 15
$ws.Range("D15").Value = '''0.6848'
$ws.Range("E15").Value = '  +0.19%  '
# Row 16
$ws.Range("D16").Value = '''83.12'
# Row 17
$ws.Range("D17").Value = '''0.000009937'
$ws.Range("E17").Value = '  +4.37%  '
# Row 18
$ws.Range("D18").Value = '''6.173'
$ws.Range("E18").Value = '  +0.48%  '
# Row 19
$ws.Range("D19").Value = '29.430.17'
$ws.Range("E19").Value = '  -0.23%  '
# Row 20
$ws.Range("D20").Value = '''231.85'
$ws.Range("E20").Value = '  -1.47%  '
# Row 21
$ws.Range("E21").Value = '  -0.35%  '
# Row 22
$ws.Range("E22").Value = '  -0.08%  '
# Row 23
$ws.Range("D23").Value = '''7.589'
$ws.Range("E23").Value = '  -1.10%  '
# Row 24
$ws.Range("D24").Value = '''0.9997'
$ws.Range("E24").Value = '  -0.09%  '
# Row 25
$ws.Range("D25").Value = '''155.10'
$ws.Range("E25").Value = '  -1.19%  '
# Row 26
$ws.Range("D26").Value = '''0.1387'
$ws.Range("E26").Value = '  -1.89%  '
# Row 27
$ws.Range("D27").Value = '''8.411'
$ws.Range("E27").Value = '  -0.80%  '
# Row 28
$ws.Range("D28").Value = '''17.66'
$ws.Range("E28").Value = '  -0.38%  '
# Row 29
$ws.Range("D29").Value = '''1.467'
$ws.Range("E29").Value = '  -1.38%  '
# Row 30
$ws.Range("D30").Value = '''0.05788'
$ws.Range("E30").Value = '  -3.39%  '
# Row 31
$ws.Range("D31").Value = '''1.259'
$ws.Range("E31").Value = '  +0.45%  '
# Row 32
$ws.Range("B32").Value = 'RocketPoolETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D32").Value = '3.692.82'
$ws.Range("E32").Value = '  +84.24%  '
# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.130'
$ws.Range("E33").Value = '  +0.22%  '
# Row 34
$ws.Range("D34").Value = '''4.024'
$ws.Range("E34").Value = '  -1.18%  '
# Row 35
$ws.Range("D35").Value = '''1.874'
$ws.Range("E35").Value = '  +0.64%  '
# Row 36
$ws.Range("D36").Value = '''1.157'
$ws.Range("E36").Value = '  -1.44%  '
# Row 37
$ws.Range("D37").Value = '''0.7169'
$ws.Range("E37").Value = '  -0.29%  '
# Row 38
$ws.Range("D38").Value = '''2.591'
$ws.Range("E38").Value = '  -0.27%  '
# Row 39
$ws.Range("D39").Value = '1.251.31'
# Row 40
$ws.Range("D40").Value = '''2.795'
$ws.Range("E40").Value = '  -0.05%  '
# Row 41
$ws.Range("D41").Value = '''0.01807'
$ws.Range("E41").Value = '  +1.76%  '
# Row 42
$ws.Range("D42").Value = '''0.9004'
$ws.Range("E42").Value = '  -1.03%  '
# Row 43
$ws.Range("D43").Value = '''6.089'
$ws.Range("E43").Value = '  -1.18%  '
# Row 44
$ws.Range("D44").Value = '''0.9994'
$ws.Range("E44").Value = '  -0.03%  '
# Row 45
$ws.Range("D45").Value = '''101.78'
$ws.Range("E45").Value = '  -0.14%  '
# Row 46
$ws.Range("E46").Value = '  +1.07%  '
# Row 47
$ws.Range("D47").Value = '''7.205'
$ws.Range("E47").Value = '  -1.44%  '
# Row 48
$ws.Range("D48").Value = '''9.175'
$ws.Range("E48").Value = '  +1.38%  '
# Row 49
$ws.Range("D49").Value = '''0.4014'
$ws.Range("E49").Value = '  -0.35%  '
# Row 50
$ws.Range("D50").Value = '''1.689'
$ws.Range("E50").Value = '  +1.99%  '
# Row 51
$ws.Range("D51").Value = '''0.1125'
$ws.Range("E51").Value = '  +0.00%  '
